$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace rows 2-20 with the new "best model" results (1 hidden layer, 70 neurons, etc.)
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = 'model_14_5_14'
$ws.Cells.Item(2, 3).Value = 0.793455378695981
$ws.Cells.Item(2, 4).Value = 0.4244238555107573
$ws.Cells.Item(2, 5).Value = 0.6608014781599317
$ws.Cells.Item(2, 6).Value = 0.5678008674674568
$ws.Cells.Item(2, 7).Value = 24.90934371948242
$ws.Cells.Item(2, 8).Value = 81.53942108154297
$ws.Cells.Item(2, 9).Value = 57.88400268554688
$ws.Cells.Item(2, 10).Value = 70.40745544433594
$ws.Cells.Item(2, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = 'model_14_5_13'
$ws.Cells.Item(3, 3).Value = 0.7933326582680003
$ws.Cells.Item(3, 4).Value = 0.4251372654216189
$ws.Cells.Item(3, 5).Value = 0.6617347785971694
$ws.Cells.Item(3, 6).Value = 0.5685893796917934
$ws.Cells.Item(3, 7).Value = 24.92414283752441
$ws.Cells.Item(3, 8).Value = 81.43834686279297
$ws.Cells.Item(3, 9).Value = 57.72473907470703
$ws.Cells.Item(3, 10).Value = 70.27899932861328
$ws.Cells.Item(3, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = 'model_14_5_15'
$ws.Cells.Item(4, 3).Value = 0.7926857422892021
$ws.Cells.Item(4, 4).Value = 0.4208282591181479
$ws.Cells.Item(4, 5).Value = 0.6596114599989675
$ws.Cells.Item(4, 6).Value = 0.5655588562216785
$ws.Cells.Item(4, 7).Value = 25.00216102600098
$ws.Cells.Item(4, 8).Value = 82.04879760742188
$ws.Cells.Item(4, 9).Value = 58.08708190917969
$ws.Cells.Item(4, 10).Value = 70.77268981933594
$ws.Cells.Item(4, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = 'model_14_5_16'
$ws.Cells.Item(5, 3).Value = 0.7926851762345462
$ws.Cells.Item(5, 4).Value = 0.4125075173457247
$ws.Cells.Item(5, 5).Value = 0.6562808351878726
$ws.Cells.Item(5, 6).Value = 0.5600862037542027
$ws.Cells.Item(5, 7).Value = 25.00222969055176
$ws.Cells.Item(5, 8).Value = 83.22755432128906
$ws.Cells.Item(5, 9).Value = 58.65545272827148
$ws.Cells.Item(5, 10).Value = 71.66421508789062
$ws.Cells.Item(5, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 2).Value = 'model_14_5_11'
$ws.Cells.Item(6, 3).Value = 0.7926035815570951
$ws.Cells.Item(6, 4).Value = 0.4327968716939861
$ws.Cells.Item(6, 5).Value = 0.6700305331880225
$ws.Cells.Item(6, 6).Value = 0.5762051655466813
$ws.Cells.Item(6, 7).Value = 25.0120677947998
$ws.Cells.Item(6, 8).Value = 80.35324859619141
$ws.Cells.Item(6, 9).Value = 56.30907440185547
$ws.Cells.Item(6, 10).Value = 69.03835296630859
$ws.Cells.Item(6, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(7, 1).Value = 12
$ws.Cells.Item(7, 2).Value = 'model_14_5_12'
$ws.Cells.Item(7, 3).Value = 0.7924749114419452
$ws.Cells.Item(7, 4).Value = 0.4250480315879348
$ws.Cells.Item(7, 5).Value = 0.6610854657107283
$ws.Cells.Item(7, 6).Value = 0.5682281144890153
$ws.Cells.Item(7, 7).Value = 25.027587890625
$ws.Cells.Item(7, 8).Value = 81.45098876953125
$ws.Cells.Item(7, 9).Value = 57.83554458618164
$ws.Cells.Item(7, 10).Value = 70.33785247802734
$ws.Cells.Item(7, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(8, 1).Value = 13
$ws.Cells.Item(8, 2).Value = 'model_14_5_17'
$ws.Cells.Item(8, 3).Value = 0.7917328975520058
$ws.Cells.Item(8, 4).Value = 0.4085525486840363
$ws.Cells.Item(8, 5).Value = 0.6530655525354701
$ws.Cells.Item(8, 6).Value = 0.5566803899427601
$ws.Cells.Item(8, 7).Value = 25.11707496643066
$ws.Cells.Item(8, 8).Value = 83.78783416748047
$ws.Cells.Item(8, 9).Value = 59.20413589477539
$ws.Cells.Item(8, 10).Value = 72.21903991699219
$ws.Cells.Item(8, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(9, 1).Value = 14
$ws.Cells.Item(9, 2).Value = 'model_14_5_10'
$ws.Cells.Item(9, 3).Value = 0.7909776258259271
$ws.Cells.Item(9, 4).Value = 0.4349555093314732
$ws.Cells.Item(9, 5).Value = 0.6747804025417153
$ws.Cells.Item(9, 6).Value = 0.5795405364507769
$ws.Cells.Item(9, 7).Value = 25.20816040039062
$ws.Cells.Item(9, 8).Value = 80.04744720458984
$ws.Cells.Item(9, 9).Value = 55.49851226806641
$ws.Cells.Item(9, 10).Value = 68.49500274658203
$ws.Cells.Item(9, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(10, 1).Value = 15
$ws.Cells.Item(10, 2).Value = 'model_14_5_18'
$ws.Cells.Item(10, 3).Value = 0.7906625088001537
$ws.Cells.Item(10, 4).Value = 0.397027494809812
$ws.Cells.Item(10, 5).Value = 0.6469965126492243
$ws.Cells.Item(10, 6).Value = 0.5483825932525763
$ws.Cells.Item(10, 7).Value = 25.24616432189941
$ws.Cells.Item(10, 8).Value = 85.42054748535156
$ws.Cells.Item(10, 9).Value = 60.23981475830078
$ws.Cells.Item(10, 10).Value = 73.57079315185547
$ws.Cells.Item(10, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(11, 1).Value = 16
$ws.Cells.Item(11, 2).Value = 'model_14_5_19'
$ws.Cells.Item(11, 3).Value = 0.7903474439790713
$ws.Cells.Item(11, 4).Value = 0.3944160967622141
$ws.Cells.Item(11, 5).Value = 0.6467230160030752
$ws.Cells.Item(11, 6).Value = 0.5470455344634841
$ws.Cells.Item(11, 7).Value = 25.28416061401367
$ws.Cells.Item(11, 8).Value = 85.79049682617188
$ws.Cells.Item(11, 9).Value = 60.28648376464844
$ws.Cells.Item(11, 10).Value = 73.78861236572266
$ws.Cells.Item(11, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(12, 1).Value = 17
$ws.Cells.Item(12, 2).Value = 'model_14_5_20'
$ws.Cells.Item(12, 3).Value = 0.7889205318091564
$ws.Cells.Item(12, 4).Value = 0.3841195382909653
$ws.Cells.Item(12, 5).Value = 0.6447578324635511
$ws.Cells.Item(12, 6).Value = 0.5413363867777254
$ws.Cells.Item(12, 7).Value = 25.45624923706055
$ws.Cells.Item(12, 8).Value = 87.24916076660156
$ws.Cells.Item(12, 9).Value = 60.62184143066406
$ws.Cells.Item(12, 10).Value = 74.71865844726562
$ws.Cells.Item(12, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(13, 1).Value = 18
$ws.Cells.Item(13, 2).Value = 'model_14_5_21'
$ws.Cells.Item(13, 3).Value = 0.7879561570013363
$ws.Cells.Item(13, 4).Value = 0.3771577554230655
$ws.Cells.Item(13, 5).Value = 0.6438289684487419
$ws.Cells.Item(13, 6).Value = 0.537673397975859
$ws.Cells.Item(13, 7).Value = 25.57254981994629
$ws.Cells.Item(13, 8).Value = 88.23540496826172
$ws.Cells.Item(13, 9).Value = 60.78035736083984
$ws.Cells.Item(13, 10).Value = 75.31536865234375
$ws.Cells.Item(13, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(14, 1).Value = 20
$ws.Cells.Item(14, 2).Value = 'model_14_5_22'
$ws.Cells.Item(14, 3).Value = 0.7871192185515063
$ws.Cells.Item(14, 4).Value = 0.3732352104046813
$ws.Cells.Item(14, 5).Value = 0.6430237927975164
$ws.Cells.Item(14, 6).Value = 0.5354705409969531
$ws.Cells.Item(14, 7).Value = 25.67348289489746
$ws.Cells.Item(14, 8).Value = 88.79109191894531
$ws.Cells.Item(14, 9).Value = 60.91775512695312
$ws.Cells.Item(14, 10).Value = 75.67423248291016
$ws.Cells.Item(14, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(15, 1).Value = 22
$ws.Cells.Item(15, 2).Value = 'model_14_5_9'
$ws.Cells.Item(15, 3).Value = 0.7849273560535748
$ws.Cells.Item(15, 4).Value = 0.42143929759103
$ws.Cells.Item(15, 5).Value = 0.6863836637931482
$ws.Cells.Item(15, 6).Value = 0.5790378269812004
$ws.Cells.Item(15, 7).Value = 25.93782424926758
$ws.Cells.Item(15, 8).Value = 81.96223449707031
$ws.Cells.Item(15, 9).Value = 53.51842498779297
$ws.Cells.Item(15, 10).Value = 68.57689666748047
$ws.Cells.Item(15, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(16, 1).Value = 23
$ws.Cells.Item(16, 2).Value = 'model_14_5_8'
$ws.Cells.Item(16, 3).Value = 0.7820880988355638
$ws.Cells.Item(16, 4).Value = 0.4299180728352087
$ws.Cells.Item(16, 5).Value = 0.6949359057733493
$ws.Cells.Item(16, 6).Value = 0.5871571571626681
$ws.Cells.Item(16, 7).Value = 26.28023910522461
$ws.Cells.Item(16, 8).Value = 80.76107788085938
$ws.Cells.Item(16, 9).Value = 52.05898666381836
$ws.Cells.Item(16, 10).Value = 67.25421142578125
$ws.Cells.Item(16, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(17, 1).Value = 27
$ws.Cells.Item(17, 2).Value = 'model_14_5_7'
$ws.Cells.Item(17, 3).Value = 0.7722166036937503
$ws.Cells.Item(17, 4).Value = 0.4088311356072766
$ws.Cells.Item(17, 5).Value = 0.7188197538103283
$ws.Cells.Item(17, 6).Value = 0.5892227943444696
$ws.Cells.Item(17, 7).Value = 27.47074317932129
$ws.Cells.Item(17, 8).Value = 83.74837493896484
$ws.Cells.Item(17, 9).Value = 47.98322296142578
$ws.Cells.Item(17, 10).Value = 66.91771697998047
$ws.Cells.Item(17, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(18, 1).Value = 29
$ws.Cells.Item(18, 2).Value = 'model_14_5_5'
$ws.Cells.Item(18, 3).Value = 0.7704571756280686
$ws.Cells.Item(18, 4).Value = 0.474431173880703
$ws.Cells.Item(18, 5).Value = 0.7483886158846897
$ws.Cells.Item(18, 6).Value = 0.6340003993584564
$ws.Cells.Item(18, 7).Value = 27.68293190002441
$ws.Cells.Item(18, 8).Value = 74.45509338378906
$ws.Cells.Item(18, 9).Value = 42.93731689453125
$ws.Cells.Item(18, 10).Value = 59.62321472167969
$ws.Cells.Item(18, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(19, 1).Value = 33
$ws.Cells.Item(19, 2).Value = 'model_14_5_6'
$ws.Cells.Item(19, 3).Value = 0.7684996735102603
$ws.Cells.Item(19, 4).Value = 0.4285630700220719
$ws.Cells.Item(19, 5).Value = 0.7326629457130962
$ws.Cells.Item(19, 6).Value = 0.6051312181445232
$ws.Cells.Item(19, 7).Value = 27.91900825500488
$ws.Cells.Item(19, 8).Value = 80.95303344726562
$ws.Cells.Item(19, 9).Value = 45.62089538574219
$ws.Cells.Item(19, 10).Value = 64.32614898681641
$ws.Cells.Item(19, 11).Value = 'Hidden Size=[70], regularizer=0.1, learning_rate=0.02'

$ws.Cells.Item(20, 1).Value = 36
$ws.Cells.Item(20, 2).Value = 'model_13_5_0'
$ws.Cells.Item(20, 3).Value = 0.7642749634634054
$ws.Cells.Item(20, 4).Value = 0.7044608375173393
$ws.Cells.Item(20, 5).Value = 0.7791132756559503
$ws.Cells.Item(20, 6).Value = 0.7271868847138594
$ws.Cells.Item(20, 7).Value = 28.42851066589355
$ws.Cells.Item(20, 8).Value = 60.08975601196289
$ws.Cells.Item(20, 9).Value = 22.10588073730469
$ws.Cells.Item(20, 10).Value = 42.21499633789062
$ws.Cells.Item(20, 11).Value = 'Hidden Size=[80], regularizer=0.02, learning_rate=0.02'

# Remove the trailing rows that are no longer part of the top-20 list
$ws.Range("A21:K24").EntireRow.Delete()
